$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BugFix 1: fixed assigning people to tickets in certain column
# Assign "Daniel" to the "Assigned To" column (B) for the rows that were
# missing an assignee.
$ws.Range("B3").Value = "Daniel"
$ws.Range("B4").Value = "Daniel"
$ws.Range("B5").Value = "Daniel"
$ws.Range("B6").Value = "Daniel"

# Move/leave the active selection where the author ended up after editing.
$ws.Range("E21").Select()
